$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new results row (row 4) to the student exam results sheet.
$ws.Range("A4").Value = "AbdulHafiz Ismail Mambo"
$ws.Range("B4").Value = "std272"
$ws.Range("C4").Value = "SS2_SILVER"
$ws.Range("D4").Value = "ACCOUNTING"

# "16%" must be stored as literal text (like the existing Score (%) cells),
# not auto-converted to a numeric percentage. Write it as a text formula,
# then convert it to a plain value in-place so no number formatting /
# formula residue is left behind on the cell.
$ws.Range("E4").Formula = '="16%"'
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial(-4163)

$ws.Range("F4").Value = 8
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = "FAIL"
$ws.Range("I4").Value = "2025-12-06 18:28"
